# Commiting Customer + Account files
# Append two new LCY Current Account rows (Customer_ID 17705019 / 17705020)
# to the bottom of the data table on Sheet0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new rows as text-returning formulas first ("=""…""") so Excel
# doesn't auto-coerce the numeric-looking IDs into Number cells - these
# columns are stored as plain text (shared strings) throughout the sheet.
$ws.Range("A17").Formula = "=""118448"""
$ws.Range("B17").Formula = "=""1008617677"""
$ws.Range("C17").Formula = "=""17705019"""
$ws.Range("D17").Formula = "=""1001"""

$ws.Range("A18").Formula = "=""118448"""
$ws.Range("B18").Formula = "=""1008617679"""
$ws.Range("C18").Formula = "=""17705020"""
$ws.Range("D18").Formula = "=""1001"""

# Convert the formulas down to plain text values (copy / paste-values) so
# the cells end up as ordinary string constants - matching the rest of the
# sheet - instead of leaving formulas behind, and without touching any
# cell styles/number formats.
$rng = $ws.Range("A17:D18")
$rng.Copy()
$rng.PasteSpecial(-4163)

Write-Host "Appended rows 17-18 to Sheet0"
